# Generate Report for Handback
#
# For each language sheet (zh-cn, de-de):
#  - Status column (B) moves from "Ready for handoff" to
#    "Handed back: in sync with en-US"
#  - New "Latest Target File" (E) / "Latest Handback File" (F) columns are
#    populated (mirroring the existing Source File (A) / Latest Handoff
#    File (C) hyperlinks, which is what a completed handback reuses)
#  - Latest Handback DateTime (G) is stamped with the handback time

$wb = $excel.ActiveWorkbook

function Set-HandbackRow {
    param($ws, $row, $statusText, $handbackDateTime)

    # Grab the existing hyperlinks for the Source File (A) and Latest
    # Handoff File (C) cells of this row so the new Latest Target File (E)
    # / Latest Handback File (F) cells can mirror them exactly.
    $aCell = $ws.Range("A" + $row)
    $cCell = $ws.Range("C" + $row)
    $aAddress = $null
    $aDisplay = $null
    $cAddress = $null
    $cDisplay = $null
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq $aCell.Address()) {
            $aAddress = $h.Address
            $aDisplay = $h.TextToDisplay
        }
        if ($h.Range.Address() -eq $cCell.Address()) {
            $cAddress = $h.Address
            $cDisplay = $h.TextToDisplay
        }
    }

    # Status text.
    $ws.Range("B" + $row).Value = $statusText

    # Latest Target File (E) mirrors the Source File hyperlink (A).
    $eCell = $ws.Range("E" + $row)
    $eCell.Value = $aDisplay
    $ws.Hyperlinks.Add($eCell, $aAddress, [System.Type]::Missing, [System.Type]::Missing, $aDisplay) | Out-Null

    # Latest Handback File (F) mirrors the Latest Handoff File hyperlink (C).
    $fCell = $ws.Range("F" + $row)
    $fCell.Value = $cDisplay
    $ws.Hyperlinks.Add($fCell, $cAddress, [System.Type]::Missing, [System.Type]::Missing, $cDisplay) | Out-Null

    # Latest Handback DateTime (G).
    $ws.Range("G" + $row).Value = $handbackDateTime
}

$zhcn = $wb.Worksheets.Item("zh-cn")
Set-HandbackRow $zhcn 2 "Handed back: in sync with en-US" "2016-02-25 04:05:26"
Set-HandbackRow $zhcn 3 "Handed back: in sync with en-US" "2016-02-25 04:05:26"

$dede = $wb.Worksheets.Item("de-de")
Set-HandbackRow $dede 2 "Handed back: in sync with en-US" "2016-02-25 04:05:53"
Set-HandbackRow $dede 3 "Handed back: in sync with en-US" "2016-02-25 04:05:53"
